$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$t.Cell(1, 1).Range.Text = "49 x 90" + $nl + "  9    0" + $nl + "  ----" + $nl + "4|    |" + $nl + "9|    |"
$t.Cell(1, 2).Range.Text = "91 x 39" + $nl + "  3    9" + $nl + "  ----" + $nl + "9|    |" + $nl + "1|    |"
$t.Cell(1, 3).Range.Text = "52 x 87" + $nl + "  8    7" + $nl + "  ----" + $nl + "5|    |" + $nl + "2|    |"
$t.Cell(2, 1).Range.Text = "93 x 83" + $nl + "  8    3" + $nl + "  ----" + $nl + "9|    |" + $nl + "3|    |"
$t.Cell(2, 2).Range.Text = "81 x 16" + $nl + "  1    6" + $nl + "  ----" + $nl + "8|    |" + $nl + "1|    |"
$t.Cell(2, 3).Range.Text = "70 x 38" + $nl + "  3    8" + $nl + "  ----" + $nl + "7|    |" + $nl + "0|    |"
$t.Cell(3, 1).Range.Text = "75 x 23" + $nl + "  2    3" + $nl + "  ----" + $nl + "7|    |" + $nl + "5|    |"
$t.Cell(3, 2).Range.Text = "20 x 66" + $nl + "  6    6" + $nl + "  ----" + $nl + "2|    |" + $nl + "0|    |"
$t.Cell(3, 3).Range.Text = "85 x 73" + $nl + "  7    3" + $nl + "  ----" + $nl + "8|    |" + $nl + "5|    |"
$t.Cell(4, 1).Range.Text = "59 x 96" + $nl + "  9    6" + $nl + "  ----" + $nl + "5|    |" + $nl + "9|    |"
$t.Cell(4, 2).Range.Text = "36 x 28" + $nl + "  2    8" + $nl + "  ----" + $nl + "3|    |" + $nl + "6|    |"
$t.Cell(4, 3).Range.Text = "12 x 19" + $nl + "  1    9" + $nl + "  ----" + $nl + "1|    |" + $nl + "2|    |"
$t.Cell(5, 1).Range.Text = "36 x 47" + $nl + "  4    7" + $nl + "  ----" + $nl + "3|    |" + $nl + "6|    |"
$t.Cell(5, 2).Range.Text = "59 x 30" + $nl + "  3    0" + $nl + "  ----" + $nl + "5|    |" + $nl + "9|    |"
$t.Cell(5, 3).Range.Text = "58 x 27" + $nl + "  2    7" + $nl + "  ----" + $nl + "5|    |" + $nl + "8|    |"
